# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
# Update the K column (column G) values on Sheet1 to the newly recalculated figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newValues = @{
    2  = 0
    3  = 4
    4  = 2
    5  = 1
    6  = 2
    7  = 1
    8  = 0
    9  = 0
    11 = 0
    12 = 4
    13 = 1
    14 = 2
    15 = 1
    16 = 1
    17 = 1
    18 = 0
    19 = 1
    20 = 0
    21 = 2
    22 = 0
    23 = 2
    24 = 1
    25 = 1
    26 = 1
    27 = 1
    28 = 0
    29 = 2
    30 = 0
    31 = 1
    32 = 0
    33 = 2
    34 = 2
    35 = 1
    36 = 2
    37 = 1
    38 = 1
    39 = 3
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
